# Applies the TEST_3_percentage_table.xlsx edit:
#  - widen column A on sheet "Section - 1" from 22 to 101 (raw OOXML width units)
#  - rename the two "P2 - Test question 2" table headers (rows 21 & 31) to
#    "P1 - Test question 1"
#  - replace the numeric "Var" codes (3/4/7/8) in the two tables below those
#    headers with their label equivalents (Blue/Green/Red/Yellow) and swap
#    around a couple of the corresponding percentage values
#  - append an explanatory note to the two related "Footer" cells (rows 28 & 38)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Section - 1")

# --- widen column A ---------------------------------------------------
# ColumnWidth (character units) and the raw OOXML <col width="..."/> differ
# by a constant padding offset (~0.8333 for this font); subtracting it here
# makes the saved worksheet XML come out to width="101".
$ws.Columns.Item(1).ColumnWidth = 101 - 0.8333333333333

# --- table headers ------------------------------------------------------
$ws.Range("A21").Value = "P1 - Test question 1"
$ws.Range("A31").Value = "P1 - Test question 1"

# --- table starting at row 21 (P4 - Example Category 1) -----------------
$ws.Range("A24").Value = "Blue"
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 0.5

$ws.Range("A25").Value = "Green"

$ws.Range("A26").Value = "Red"

$ws.Range("A27").Value = "Yellow"
$ws.Range("C27").Value = 0.6666666666666666
$ws.Range("D27").Value = 0.3333333333333333

$ws.Range("A28").Value = "Footer // This table was calculated by a variable operation P2 - Test question 2 * Expansion_factor"

# --- table starting at row 31 (P5 - Example Category 2) -----------------
$ws.Range("A34").Value = "Blue"
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 0.5
$ws.Range("E34").Value = 0.5

$ws.Range("A35").Value = "Green"
$ws.Range("C35").Value = 0.5
$ws.Range("E35").Value = 0

$ws.Range("A36").Value = "Red"

$ws.Range("A37").Value = "Yellow"
$ws.Range("C37").Value = 0.3333333333333333
$ws.Range("D37").Value = 0.3333333333333333
$ws.Range("E37").Value = 0.3333333333333333

$ws.Range("A38").Value = "Footer // This table was calculated by a variable operation P2 - Test question 2 * Expansion_factor"
